$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the value from B13 (DrugUse.ObservationOfUse.StartDate) down to B16,
# matching the formatting of the rest of column B, then fully clear the
# now-empty B13 cell and update the active selection to match.
$ws.Range("B16").Value = $ws.Range("B13").Value2
$ws.Range("B16").Font.Name = "Calibri"
$ws.Range("B13").Clear()

$ws.Range("B16").Select()
